# S9 Table update: refresh gene counts / percentages with revised data
# (121 -> 131 genes), re-sort the metabolic-process table by the new
# frequency counts, and switch the percentage column to a 1-decimal
# number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Title text: "121 genes" -> "131 genes" (A1), keeping the italic
#    "Bordetella pertussis" run intact.
# ---------------------------------------------------------------------
$titleCell = $ws.Range("A1")
$fullText = $titleCell.Text
$needle = "121 genes"
$idx = $fullText.IndexOf($needle)
if ($idx -ge 0) {
    $titleCell.Characters($idx + 1, 3).Text = "131"
}

$fullText = $titleCell.Text
$species = "Bordetella pertussis"
$speciesIdx = $fullText.IndexOf($species)
if ($speciesIdx -ge 0) {
    $titleCell.Characters($speciesIdx + 1, $species.Length).Font.Italic = $true
}

# ---------------------------------------------------------------------
# 2. Refresh the data table (rows 4-19) with the revised counts /
#    percentages, re-sorted in descending frequency order.
# ---------------------------------------------------------------------
$data = @(
    @("Carbohydrate metabolism", 75, 28.846153846153847),
    @("Amino acid metabolism", 59, 22.692307692307693),
    @("Metabolism of cofactors and vitamins", 26, 10),
    @("Nucleotide metabolism", 18, 6.9230769230769234),
    @("Energy metabolism", 17, 6.5384615384615383),
    @("Lipid metabolism", 14, 5.384615384615385),
    @("Metabolism of other amino acids", 13, 5),
    @("Glycan biosynthesis and metabolism", 10, 3.8461538461538463),
    @("Metabolism of terpenoids and polyketides", 8, 3.0769230769230771),
    @("Biosynthesis of other secondary metabolites", 5, 1.9230769230769231),
    @("Quorum sensing", 4, 1.5384615384615385),
    @("Xenobiotics biodegradation and metabolism", 4, 1.5384615384615385),
    @("Vancomycin resistance", 3, 1.1538461538461537),
    @("Cationic antimicrobial peptide (CAMP) resistance", 2, 0.76923076923076927),
    @("Two-component system", 1, 0.38461538461538464),
    @("Aminoacyl-tRNA biosynthesis", 1, 0.38461538461538464)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# ---------------------------------------------------------------------
# 3. Formatting: columns A/B for the data rows lose their old "center"
#    style (back to workbook default); column C keeps a numeric style
#    but now formatted with a single decimal place instead of two.
# ---------------------------------------------------------------------
$ws.Range("A4:B19").ClearFormats()
$ws.Range("C4:C19").ClearFormats()
$ws.Range("C4:C19").NumberFormat = "0.0"

# ---------------------------------------------------------------------
# 4. Update the view: the sheet is scrolled down with C4:C19 selected.
# ---------------------------------------------------------------------
$ws.Range("C4:C19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4

Write-Host "S9 table refreshed"
